$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price cells remain stored as text (matching source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "60.858.63"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "3.364.93"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "571.46"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "135.65"
$ws.Range("E6").Value = "  +8.16%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.364.96"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D10").Value = "7.59"
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").Value = "3.936.69"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "3.365.38"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "25.18"
$ws.Range("E17").Value = "  +3.24%  "
$ws.Range("D18").Value = "61.003.47"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").Value = "13.90"
$ws.Range("E19").Value = "  +6.33%  "
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("D22").Value = "372.48"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "0.568"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("D24").Value = "3.499.33"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "70.60"
$ws.Range("E26").Value = "  -1.39%  "
$ws.Range("E27").Value = "  +10.72%  "
$ws.Range("E28").Value = "  +22.26%  "
$ws.Range("D29").Value = "7.69"
$ws.Range("E29").Value = "  +11.17%  "
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").Value = "0.155"
$ws.Range("E33").Value = "  +4.94%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "3.395.79"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").Value = "23.37"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "5.57"
$ws.Range("E37").Value = "  +6.22%  "
$ws.Range("E38").Value = "  +4.23%  "
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("D40").Value = "163.15"
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "0.0785"
$ws.Range("E41").Value = "  +4.29%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("E43").Value = "  +3.94%  "
$ws.Range("E44").Value = "  +11.93%  "
$ws.Range("E45").Value = "  -0.54%  "
$ws.Range("D46").Value = "41.31"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").Value = "23.02"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").Value = "6.97"
$ws.Range("E49").Value = "  +5.92%  "
$ws.Range("D50").Value = "23.23"
$ws.Range("E50").Value = "  +15.47%  "
$ws.Range("E51").Value = "  +13.85%  "
